$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-19 Saturday", "2023-08-20 Sunday"),
    @("86÷7=", "72÷3="),
    @("10÷7=", "20÷6="),
    @("90÷3=", "42÷5="),
    @("56÷8=", "55÷9="),
    @("74÷3=", "96÷4="),
    @("17÷6=", "60÷7="),
    @("59÷7=", "57÷5="),
    @("56÷9=", "26÷6="),
    @("62÷8=", "78÷8="),
    @("99÷7=", "87÷2="),
    @("83÷6=", "84÷9="),
    @("40÷7=", "92÷4="),
    @("44÷6=", "13÷5="),
    @("79÷9=", "79÷4="),
    @("75÷9=", "44÷9="),
    @("71÷3=", "57÷9="),
    @("77÷2=", "52÷7="),
    @("72÷4=", "41÷7="),
    @("89÷9=", "93÷5="),
    @("59÷3=", "42÷8="),
    @("33÷4=", "21÷7="),
    @("29÷3=", "63÷9="),
    @("59÷9=", "44÷5="),
    @("25÷6=", "95÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
